$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.986
$ws.Range("D2").Value = 0.9859838095952398
$ws.Range("E2").Value = 0.998
$ws.Range("F2").Value = 0.996
$ws.Range("G2").Value = 0.9728643216080402
$ws.Range("H2").Value = 0.974
$ws.Range("I2").Value = 0.9890547263681593

# Row 3
$ws.Range("C3").Value = 0.9928
$ws.Range("D3").Value = 0.9927939867747483
$ws.Range("F3").Value = 0.9990009990009989
$ws.Range("H3").Value = 0.986
$ws.Range("I3").Value = 0.996

# Row 4
$ws.Range("C4").Value = 0.994
$ws.Range("D4").Value = 0.993998571303306
$ws.Range("G4").Value = 0.9860834990059641
$ws.Range("H4").Value = 0.9869083585095669

# Row 5
$ws.Range("C5").Value = 0.9944
$ws.Range("D5").Value = 0.9943983379159729
$ws.Range("G5").Value = 0.9870388833499502
$ws.Range("H5").Value = 0.9879518072289156

# Row 6
$ws.Range("C6").Value = 0.9944
$ws.Range("D6").Value = 0.9943975967975968
$ws.Range("E6").Value = 0.9990009990009989
$ws.Range("G6").Value = 0.988
$ws.Range("H6").Value = 0.988988988988989
$ws.Range("I6").Value = 0.9969969969969971
